# Update the "( N – 4 )" verse-range captions to use "/" instead of the
# en-dash "–" as the separator, e.g. "( 3 – 4 )" -> "( 3 / 4 )".
#
# Each affected slide has two shapes:
#   Shape 1 = the lyric content placeholder (untouched)
#   Shape 2 = "TextBox 2" holding the caption, e.g. "( 3 – 4 )"
#
# Two different source layouts exist for Shape 2:
#   (a) a single run containing the whole "( N – 4 )" string
#       (slides 3, 7, 11, 15 -- titles)
#   (b) three runs "( " / "N " / "– 4 )"
#       (slides 6, 10, 14 -- subtitles), or "( " / "N – 4 " / ")"
#       (slide 2 -- subtitle)
#
# For layout (a) we simply overwrite the whole TextRange.Text, which keeps
# the single run's formatting intact.
#
# For layout (b) we need the dash run to be split into a new "/ " run and
# the remainder, with the new run picking up the same formatting
# (including smtClean="0") as its neighbouring "N " run. Doing the
# replacement in two steps (first swap in " / " across the run boundary,
# then re-glue the leading digit + space back onto the first run) makes
# the host attribute the new run's formatting to the digit run, which is
# what the target markup has.

$p = $ppt.ActivePresentation

function Update-SingleRunCaption($slideIndex, $newText) {
    $s = $p.Slides.Item($slideIndex)
    $shp = $s.Shapes.Item(2)
    $shp.TextFrame.TextRange.Text = $newText
}

function Update-SplitRunCaption($slideIndex, $digit) {
    $s = $p.Slides.Item($slideIndex)
    $shp = $s.Shapes.Item(2)
    $tr = $shp.TextFrame.TextRange
    # Text is "( N - 4 )" (9 chars): swap " - " (chars 4-6) for " / ".
    $dash = $tr.Characters(4, 3)
    $dash.Text = " / "
    # Re-attach the digit + trailing space ("N ", chars 3-4) as one run so
    # the new "/ " run's rPr is inherited from the digit run (smtClean="0").
    $digitRun = $tr.Characters(3, 2)
    $digitRun.Text = "$digit "
}

# Subtitle slides with 3 runs: "( " / "N " / "- 4 )"
Update-SplitRunCaption 10 "3"
Update-SplitRunCaption 14 "4"
Update-SplitRunCaption 6  "2"

# Title slides with a single run: "( N - 4 )"
Update-SingleRunCaption 11 "( 3 / 4 )"
Update-SingleRunCaption 15 "( 4 / 4 )"
Update-SingleRunCaption 3  "( 1 / 4 )"
Update-SingleRunCaption 7  "( 2 / 4 )"

# Slide 2 subtitle has runs "( " / "1 - 4 " / ")": swap the inner "- " for
# "/ ", splitting the middle run into "1 " / "/ " / "4 " (all three keep
# smtClean="0", inherited from the original middle run).
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange
$dash2 = $tr2.Characters(5, 2)
$dash2.Text = "/ "
